$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 20333.334
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -770

# Row 32
$ws.Range("H32").Value = 27563.125
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 27563.125
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 27563.125
$ws.Range("N32").Value = -28215.125
$ws.Range("M32").ClearContents()

# Row 74
$ws.Range("H74").Value = 16745507
$ws.Range("I74").Value = 25115824
$ws.Range("J74").Value = 4875
$ws.Range("K74").Value = 25115824
$ws.Range("L74").Value = 4875
$ws.Range("M74").Value = -25114888
$ws.Range("N74").Value = -6747

# Row 77
$ws.Range("H77").Value = 16745507
$ws.Range("I77").Value = 25115824
$ws.Range("J77").Value = 4875
$ws.Range("K77").Value = 125579120
$ws.Range("L77").Value = 24375
$ws.Range("M77").Value = -125574440
$ws.Range("N77").Value = -33735

# Row 98
$ws.Range("H98").Value = 5287.8696
$ws.Range("I98").Value = 3255.5
$ws.Range("J98").Value = 50000
$ws.Range("K98").Value = 3255.5
$ws.Range("L98").Value = 50000
$ws.Range("M98").Value = -1757.5
$ws.Range("N98").Value = -52996

# Row 113
$ws.Range("H113").Value = 2861
$ws.Range("J113").Value = 3166.6667
$ws.Range("L113").Value = 3166.6667
$ws.Range("N113").Value = -9674.6667

# Row 122
$ws.Range("H122").Value = 5287.8696
$ws.Range("I122").Value = 3255.5
$ws.Range("J122").Value = 50000
$ws.Range("K122").Value = 9766.5
$ws.Range("L122").Value = 150000
$ws.Range("M122").Value = -7316.5
$ws.Range("N122").Value = -154900


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 111126000
$ws.Range("J43").Value = 111126000
$ws.Range("L43").Value = 111126000
$ws.Range("N43").Value = -111126626

# Row 74
$ws.Range("H74").Value = 1770
$ws.Range("I74").Value = 1837.5
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1837.5
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -963.5
$ws.Range("N74").Value = -3248

# Row 77
$ws.Range("H77").Value = 1770
$ws.Range("I77").Value = 1837.5
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 9187.5
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -4819.5
$ws.Range("N77").Value = -16236

# Row 132
$ws.Range("H132").Value = 14111.454
$ws.Range("I132").Value = 30628
$ws.Range("J132").Value = 4673.4287
$ws.Range("K132").Value = 91884
$ws.Range("L132").Value = 14020.2861
$ws.Range("M132").Value = -89354
$ws.Range("N132").Value = -19080.2861


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 61465
$ws.Range("I86").Value = 2888.875
$ws.Range("J86").Value = 113532.664
$ws.Range("K86").Value = 2888.875
$ws.Range("L86").Value = 113532.664
$ws.Range("M86").Value = -1765.875
$ws.Range("N86").Value = -115778.664

# Row 89
$ws.Range("H89").Value = 61465
$ws.Range("I89").Value = 2888.875
$ws.Range("J89").Value = 113532.664
$ws.Range("K89").Value = 14444.375
$ws.Range("L89").Value = 567663.3200000001
$ws.Range("M89").Value = -8828.375
$ws.Range("N89").Value = -578895.3200000001

# Row 134
$ws.Range("H134").Value = 2961.7742
$ws.Range("I134").Value = 2770.75
$ws.Range("J134").Value = 3309.0908
$ws.Range("K134").Value = 8312.25
$ws.Range("L134").Value = 9927.2724
$ws.Range("M134").Value = -5777.25
$ws.Range("N134").Value = -14997.2724


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 370.75
$ws.Range("J22").Value = 685
$ws.Range("L22").Value = 685
$ws.Range("N22").Value = -1385

# Row 62
$ws.Range("H62").Value = 35953.668
$ws.Range("I62").Value = 44108.75
$ws.Range("J62").Value = 3333.3333
$ws.Range("K62").Value = 44108.75
$ws.Range("L62").Value = 3333.3333
$ws.Range("M62").Value = -43484.75
$ws.Range("N62").Value = -4581.3333

# Row 65
$ws.Range("H65").Value = 35953.668
$ws.Range("I65").Value = 44108.75
$ws.Range("J65").Value = 3333.3333
$ws.Range("K65").Value = 220543.75
$ws.Range("L65").Value = 16666.6665
$ws.Range("M65").Value = -217423.75
$ws.Range("N65").Value = -22906.6665

# Row 80
$ws.Range("H80").Value = 34128
$ws.Range("J80").Value = 34128
$ws.Range("L80").Value = 34128
$ws.Range("N80").Value = -36374

# Row 83
$ws.Range("H83").Value = 34128
$ws.Range("J83").Value = 34128
$ws.Range("L83").Value = 102384
$ws.Range("N83").Value = -113616


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 7253.6924
$ws.Range("J80").Value = 3474.75
$ws.Range("L80").Value = 10424.25
$ws.Range("N80").Value = -12296.25

# Row 83
$ws.Range("H83").Value = 7253.6924
$ws.Range("J83").Value = 3474.75
$ws.Range("L83").Value = 31272.75
$ws.Range("N83").Value = -40632.75


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 4941.8
$ws.Range("J49").Value = 4941.8
$ws.Range("L49").Value = 4941.8
$ws.Range("N49").Value = -5309.8

# Row 97
$ws.Range("H97").Value = 36926.355
$ws.Range("I97").Value = 39458.31
$ws.Range("J97").Value = 4011
$ws.Range("K97").Value = 39458.31
$ws.Range("L97").Value = 4011
$ws.Range("M97").Value = -38962.31
$ws.Range("N97").Value = -5003


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# Row 42
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 49
$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 61
$ws.Range("H61").Value = 2776.7058
$ws.Range("I61").Value = 2738.7693
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 2738.7693
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -2536.7693
$ws.Range("N61").Value = -3304

# Row 113
$ws.Range("H113").Value = 2776.7058
$ws.Range("I113").Value = 2738.7693
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 2738.7693
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -568.7692999999999
$ws.Range("N113").Value = -7240


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 5025000
$ws.Range("I20").Value = 5025000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 5025000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -5024760
$ws.Range("N20").ClearContents()

# Row 41
$ws.Range("H41").Value = 333341400
$ws.Range("J41").Value = 333341400
$ws.Range("L41").Value = 333341400
$ws.Range("N41").Value = -333342180

# Row 94
$ws.Range("H94").Value = 49266.668
$ws.Range("J94").Value = 49266.668
$ws.Range("L94").Value = 49266.668
$ws.Range("N94").Value = -51068.668

# Row 125
$ws.Range("H125").Value = 59079.5
$ws.Range("J125").Value = 59079.5
$ws.Range("L125").Value = 59079.5

